# eadmin report and config
# Rename the "state" field to "district" and swap the state-related sample
# values/lists for Tamil Nadu district data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: D1 "state" -> "district"
$ws.Range("D1").Value = "district"

# Row 2: full list of all states -> full list of all TN districts
$ws.Range("D2").Value = "ARIYALUR,CHENGALPATTU,CHENNAI,COIMBATORE,CUDDALORE,DHARMAPURI,DINDIGUL,ERODE,KALLAKURICHI,KANCHEEPURAM,KANNIYAKUMARI,KARUR,KRISHNAGIRI,MADURAI,MAYILADUTHURAI,NAGAPATTINAM,NAMAKKAL,PERAMBALUR,PUDUKKOTTAI,RAMANATHAPURAM,RANIPET,SALEM,SIVAGANGA,TENKASI,THANJAVUR,THE NILGIRIS,THENI,THIRUVALLUR,THIRUVARUR,THOOTHUKKUDI,TIRUCHIRAPPALLI,TIRUNELVELI,TIRUPATHUR,TIRUPPUR,TIRUVANNAMALAI,VELLORE,VILUPPURAM,VIRUDHUNAGAR"

# Row 3: sub-list of states -> single selected district
$ws.Range("D3").Value = "CHENGALPATTU"

# Move the active selection from F3 to D5
$ws.Range("D5").Select()
